$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 13:22"

# Row 12 - Iran
$ws.Range("B12").Value = 92584
$ws.Range("C12").Value = 1112
$ws.Range("D12").Value = 72439
$ws.Range("E12").Value = 14268
$ws.Range("F12").Value = 2983
$ws.Range("G12").Value = 71
$ws.Range("H12").Value = 5877

# Row 60 - Kuwait
$ws.Range("B60").Value = 3440
$ws.Range("C60").Value = 152
$ws.Range("D60").Value = 2241
$ws.Range("E60").Value = 1176
$ws.Range("F60").Value = 67
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 23

# Row 63 - Barein
$ws.Range("B63").Value = 2810
$ws.Range("C63").Value = 87
$ws.Range("D63").Value = 1246
$ws.Range("E63").Value = 1556

# Row 68 - Uzbekistan
$ws.Range("D68").Value = 958
$ws.Range("E68").Value = 973

# Row 75 - Estonia
$ws.Range("B75").Value = 1660
$ws.Range("C75").Value = 13
$ws.Range("D75").Value = 240
$ws.Range("E75").Value = 1370
$ws.Range("F75").Value = 9

# Row 156 - Uganda
$ws.Range("D156").Value = 52
$ws.Range("E156").Value = 27
